$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- Paragraph 1: "DataGenerator" heading -> drop the spellStart/spellEnd proofErr markers ---
$p1 = $d.Paragraphs(1)
$xml1 = "<w:p xmlns:w='$wNs'><w:pPr><w:pStyle w:val='berschrift1'/><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:lang w:val='de-DE'/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:lang w:val='de-DE'/></w:rPr><w:t>DataGenerator</w:t></w:r></w:p>"
$p1.Range.InsertXML($xml1)

# --- Paragraph 3: "Liniendicke " + "= ?" -> single run "Liniendicke = ?", drop gramStart/gramEnd ---
$p3 = $d.Paragraphs(3)
$xml3 = "<w:p xmlns:w='$wNs'><w:pPr><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:lang w:val='de-DE'/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:lang w:val='de-DE'/></w:rPr><w:t>Liniendicke = ?</w:t></w:r></w:p>"
$p3.Range.InsertXML($xml3)

# --- Paragraph 4: long "Wenn liniendicke ..." paragraph -> merge all runs into one, drop all proofErr ---
$p4 = $d.Paragraphs(4)
$text4 = "Wenn liniendicke = 7 (wie bei angleNetwork): Wenn ich will dass immer die volle Breite des Striches im Bild ist, dann ist der Bereich in dem der Strich sich bewegen kann bei imageSize = (29, 29)  =&gt;"
$xml4 = "<w:p xmlns:w='$wNs'><w:pPr><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:lang w:val='de-DE'/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:lang w:val='de-DE'/></w:rPr><w:t>$text4</w:t></w:r></w:p>"
$p4.Range.InsertXML($xml4)

# --- Paragraph 6: "Entschieden:" -> append new run " 7" ---
$p6 = $d.Paragraphs(6)
$xml6 = "<w:p xmlns:w='$wNs'><w:pPr><w:pStyle w:val='berschrift2'/><w:rPr><w:lang w:val='de-DE'/></w:rPr></w:pPr><w:r><w:t>Entschieden</w:t></w:r><w:r><w:rPr><w:lang w:val='de-DE'/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:lang w:val='de-DE'/></w:rPr><w:t xml:space='preserve'> 7</w:t></w:r></w:p>"
$p6.Range.InsertXML($xml6)

# --- After paragraph 7 (first empty paragraph following "Entschieden:"), insert two new paragraphs ---
$p7 = $d.Paragraphs(7)
$p7.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs(8)
$p8.Range.InsertParagraphAfter()

$p8 = $d.Paragraphs(8)
$xml8 = "<w:p xmlns:w='$wNs'><w:pPr><w:rPr><w:lang w:val='de-DE'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='de-DE'/></w:rPr><w:t>Wie mache ich es, dass prior neuronen mehr impact haben?</w:t></w:r></w:p>"
$p8.Range.InsertXML($xml8)

$p9 = $d.Paragraphs(9)
$xml9 = "<w:p xmlns:w='$wNs'><w:pPr><w:rPr><w:lang w:val='de-DE'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='de-DE'/></w:rPr><w:t>Entweder mehr A Neuronen machen, oder ATilde größer machen?</w:t></w:r></w:p>"
$p9.Range.InsertXML($xml9)
